$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")
$c = $ws.Range("E2")
Write-Host "E2 top border:" $c.Borders.Item(8).LineStyle
Write-Host "E2 bottom border:" $c.Borders.Item(9).LineStyle
Write-Host "E2 left border:" $c.Borders.Item(7).LineStyle
Write-Host "E2 right border:" $c.Borders.Item(10).LineStyle
$c2 = $ws.Range("A2")
Write-Host "A2 (s2) top border:" $c2.Borders.Item(8).LineStyle
